$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7854210155738883
$ws.Range("C2").Value = 0.1977490234365575
$ws.Range("E2").Value = 0.5377489462377554
$ws.Range("F2").Value = 1.886021615883919
$ws.Range("G2").Value = 0.002400000844395342
$ws.Range("I2").Value = 0.363236892907203
$ws.Range("J2").Value = 0.03066325496517486
$ws.Range("M2").Value = 0.5869644130803451
$ws.Range("O2").Value = 1.663272803194701

$ws.Range("B3").Value = 0.6874601688413122
$ws.Range("C3").Value = 0.1735870702608224
$ws.Range("E3").Value = 0.5332302075385584
$ws.Range("F3").Value = 1.881149085181249
$ws.Range("G3").Value = 0.00240261849205622
$ws.Range("I3").Value = 0.3718039240765876
$ws.Range("J3").Value = 0.03093978266645969
$ws.Range("M3").Value = 0.5469028548780273
$ws.Range("O3").Value = 1.683161982392434

$ws.Range("B4").Value = 0.6271188673048016
$ws.Range("C4").Value = 0.1586822881201613
$ws.Range("E4").Value = 0.5306843695216656
$ws.Range("F4").Value = 1.879414442267731
$ws.Range("G4").Value = 0.002404310250795848
$ws.Range("I4").Value = 0.3774169532698632
$ws.Range("J4").Value = 0.03112717696169298
$ws.Range("M4").Value = 0.5224072806424331
$ws.Range("O4").Value = 1.697001331041619

$ws.Range("B5").Value = 0.6024823032914242
$ws.Range("C5").Value = 0.1525914082559439
$ws.Range("E5").Value = 0.5297044750455271
$ws.Range("F5").Value = 1.879023290837125
$ws.Range("G5").Value = 0.002405020970864648
$ws.Range("I5").Value = 0.3797928026534532
$ws.Range("J5").Value = 0.03120797312889678
$ws.Range("M5").Value = 0.5124514094251111
$ws.Range("O5").Value = 1.703049029124827

$ws.Range("B6").Value = 0.5983886274744066
$ws.Range("C6").Value = 0.1515790023475176
$ws.Range("E6").Value = 0.5295452420132847
$ws.Range("F6").Value = 1.878977398621927
$ws.Range("G6").Value = 0.002405140274490269
$ws.Range("I6").Value = 0.3801926487487091
$ws.Range("J6").Value = 0.03122165705605795
$ws.Range("M6").Value = 0.5107998471607971
$ws.Range("O6").Value = 1.704077856331224

$ws.Range("B7").Value = 0.6267867982304551
$ws.Range("C7").Value = 0.1586002129895121
$ws.Range("E7").Value = 0.5306709212165259
$ws.Range("F7").Value = 1.879407889194084
$ws.Range("G7").Value = 0.002404319749195965
$ws.Range("I7").Value = 0.37744863680037
$ws.Range("J7").Value = 0.031128248658014
$ws.Range("M7").Value = 0.5222729051938018
$ws.Range("O7").Value = 1.697081241760429

$ws.Range("B8").Value = 0.7516850416608918
$ws.Range("C8").Value = 0.189432556841922
$ws.Range("E8").Value = 0.5361434487158405
$ws.Range("F8").Value = 1.884080417294427
$ws.Range("G8").Value = 0.002400885907424781
$ws.Range("I8").Value = 0.3661174370099491
$ws.Range("J8").Value = 0.03075495041217025
$ws.Range("M8").Value = 0.5731302408728496
$ws.Range("O8").Value = 1.669792152299792

$ws.Range("B9").Value = 0.9950244096683036
$ws.Range("C9").Value = 0.2493327391674995
$ws.Range("E9").Value = 0.5486885322699209
$ws.Range("F9").Value = 1.90323845412135
$ws.Range("G9").Value = 0.002394819798653006
$ws.Range("I9").Value = 0.3467065506732681
$ws.Range("J9").Value = 0.03016242044265027
$ws.Range("M9").Value = 0.6736558225265981
$ws.Range("O9").Value = 1.629242619880529

$ws.Range("B10").Value = 1.172783483433477
$ws.Range("C10").Value = 0.292986341311348
$ws.Range("E10").Value = 0.5590112674147107
$ws.Range("F10").Value = 1.923440706353176
$ws.Range("G10").Value = 0.002390765937278728
$ws.Range("I10").Value = 0.3341713694016111
$ws.Range("J10").Value = 0.0298119137457693
$ws.Range("M10").Value = 0.7479801970523567
$ws.Range("O10").Value = 1.607427270509447

$ws.Range("B11").Value = 1.253418219181299
$ws.Range("C11").Value = 0.3127660710422617
$ws.Range("E11").Value = 0.5639476725198591
$ws.Range("F11").Value = 1.93396915253642
$ws.Range("G11").Value = 0.00238900834378586
$ws.Range("I11").Value = 0.328846387818766
$ws.Range("J11").Value = 0.02967083560602113
$ws.Range("M11").Value = 0.7818910902111185
$ws.Range("O11").Value = 1.599250615895016

$ws.Range("B12").Value = 1.283918353165177
$ws.Range("C12").Value = 0.3202445648337857
$ws.Range("E12").Value = 0.5658515331593179
$ws.Range("F12").Value = 1.938148968982119
$ws.Range("G12").Value = 0.002388355165659606
$ws.Range("I12").Value = 0.3268844609347354
$ws.Range("J12").Value = 0.02962005115759858
$ws.Range("M12").Value = 0.7947462874113569
$ws.Range("O12").Value = 1.596406849937779

$ws.Range("B13").Value = 1.277351162788023
$ws.Range("C13").Value = 0.3186344618534918
$ws.Range("E13").Value = 0.5654399666344858
$ws.Range("F13").Value = 1.937240182952209
$ws.Range("G13").Value = 0.002388495289387204
$ws.Range("I13").Value = 0.3273045675807982
$ws.Range("J13").Value = 0.02963087117769447
$ws.Range("M13").Value = 0.791977082231341
$ws.Range("O13").Value = 1.597008052637079

$ws.Range("B14").Value = 1.25592818387895
$ws.Range("C14").Value = 0.3133815670313709
$ws.Range("E14").Value = 0.5641036120400287
$ws.Range("F14").Value = 1.934309159402858
$ws.Range("G14").Value = 0.002388954358619209
$ws.Range("I14").Value = 0.3286838845566749
$ws.Range("J14").Value = 0.02966660466862692
$ws.Range("M14").Value = 0.7829484200437662
$ws.Range("O14").Value = 1.599011588520028

$ws.Range("B15").Value = 1.24280145198685
$ws.Range("C15").Value = 0.3101624877966174
$ws.Range("E15").Value = 0.5632895547234398
$ws.Range("F15").Value = 1.93253896100623
$ws.Range("G15").Value = 0.002389237162802769
$ws.Range("I15").Value = 0.3295358651383378
$ws.Range("J15").Value = 0.0296888360340084
$ws.Range("M15").Value = 0.7774198947108317
$ws.Range("O15").Value = 1.60027174089177

$ws.Range("B16").Value = 1.167509073862391
$ws.Range("C16").Value = 0.2916920778103531
$ws.Range("E16").Value = 0.5586934987330494
$ws.Range("F16").Value = 1.922779618693653
$ws.Range("G16").Value = 0.002390882536356323
$ws.Range("I16").Value = 0.3345269853118396
$ws.Range("J16").Value = 0.02982150293528285
$ws.Range("M16").Value = 0.7457660125070902
$ws.Range("O16").Value = 1.607996904800899

$ws.Range("B17").Value = 1.121259907698743
$ws.Range("C17").Value = 0.2803407056125593
$ws.Range("E17").Value = 0.5559355480033119
$ws.Range("F17").Value = 1.917135694205953
$ws.Range("G17").Value = 0.002391914040095484
$ws.Range("I17").Value = 0.3376857110462321
$ws.Range("J17").Value = 0.02990759260106124
$ws.Range("M17").Value = 0.7263726904885743
$ws.Range("O17").Value = 1.613184484851303

$ws.Range("B18").Value = 1.094637173025262
$ws.Range("C18").Value = 0.2738043284558671
$ws.Range("E18").Value = 0.5543718907883317
$ws.Range("F18").Value = 1.914015400887479
$ws.Range("G18").Value = 0.002392515481547443
$ws.Range("I18").Value = 0.3395380271279187
$ws.Range("J18").Value = 0.02995883823679613
$ws.Range("M18").Value = 0.7152276439951493
$ws.Range("O18").Value = 1.616332626175847

$ws.Range("B19").Value = 1.085619539470429
$ws.Range("C19").Value = 0.2715899693367305
$ws.Range("E19").Value = 0.5538463534823208
$ws.Range("F19").Value = 1.912980539643485
$ws.Range("G19").Value = 0.00239272052045969
$ws.Range("I19").Value = 0.3401712786264106
$ws.Range("J19").Value = 0.02997648619424353
$ws.Range("M19").Value = 0.7114557678084452
$ws.Range("O19").Value = 1.617426725857001

$ws.Range("B20").Value = 1.126185440704717
$ws.Range("C20").Value = 0.2815498438298221
$ws.Range("E20").Value = 0.5562267933713514
$ws.Range("F20").Value = 1.917723461585837
$ws.Range("G20").Value = 0.002391803391635793
$ws.Range("I20").Value = 0.3373457831808579
$ws.Range("J20").Value = 0.02989824926960516
$ws.Range("M20").Value = 0.7284361656723064
$ws.Range("O20").Value = 1.612615235678447

$ws.Range("B21").Value = 1.262221580776327
$ws.Range("C21").Value = 0.314924788970302
$ws.Range("E21").Value = 0.5644951943705649
$ws.Range("F21").Value = 1.935164832458568
$ws.Range("G21").Value = 0.002388819183272724
$ws.Range("I21").Value = 0.3282772637167604
$ws.Range("J21").Value = 0.0296560372809882
$ws.Range("M21").Value = 0.7855999847351001
$ws.Range("O21").Value = 1.598416236815325

$ws.Range("B22").Value = 1.350927297074634
$ws.Range("C22").Value = 0.3366691090267011
$ws.Range("E22").Value = 0.5701004458963439
$ws.Range("F22").Value = 1.947688494386142
$ws.Range("G22").Value = 0.002386940988857278
$ws.Range("I22").Value = 0.3226684174592744
$ws.Range("J22").Value = 0.02951311745164276
$ws.Range("M22").Value = 0.8230404824900432
$ws.Range("O22").Value = 1.590609059895087

$ws.Range("B23").Value = 1.303602371150646
$ws.Range("C23").Value = 0.3250701120235249
$ws.Range("E23").Value = 0.5670904044013554
$ws.Range("F23").Value = 1.940901313005796
$ws.Range("G23").Value = 0.002387936833079062
$ws.Range("I23").Value = 0.325632782916685
$ws.Range("J23").Value = 0.02958799005369528
$ws.Range("M23").Value = 0.803050587685334
$ws.Range("O23").Value = 1.59464071660517

$ws.Range("B24").Value = 1.123958711379714
$ws.Range("C24").Value = 0.2810032245674563
$ws.Range("E24").Value = 0.5560950530411191
$ws.Range("F24").Value = 1.917457344256192
$ws.Range("G24").Value = 0.002391853389673404
$ws.Range("I24").Value = 0.3374993514448104
$ws.Range("J24").Value = 0.02990246793319251
$ws.Range("M24").Value = 0.7275032546943834
$ws.Range("O24").Value = 1.612872077074712

$ws.Range("B25").Value = 0.9293702447085366
$ws.Range("C25").Value = 0.233189619263527
$ws.Range("E25").Value = 0.5451005804762161
$ws.Range("F25").Value = 1.896982179821393
$ws.Range("G25").Value = 0.002396389791790306
$ws.Range("I25").Value = 0.3516555606254066
$ws.Range("J25").Value = 0.03030780390111687
$ws.Range("M25").Value = 0.6463776178138545
$ws.Range("O25").Value = 1.638816867142381
